# Scheduled runner: refresh market-price-derived Leve profit figures
# across the per-job Sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Values come from currentAveragePrice(NQ/HQ) market syncs; columns:
#   H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#   K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 379.7
$ws.Range("I2").Value = 339.5
$ws.Range("J2").Value = 440
$ws.Range("K2").Value = 339.5
$ws.Range("L2").Value = 440
$ws.Range("M2").Value = -226.5
$ws.Range("N2").Value = -666

# Row 39
$ws.Range("H39").Value = 157.7
$ws.Range("I39").Value = 176
$ws.Range("J39").Value = 84.5
$ws.Range("K39").Value = 528
$ws.Range("L39").Value = 253.5
$ws.Range("M39").Value = -232
$ws.Range("N39").Value = -845.5

# Row 58
$ws.Range("H58").Value = 967.9231
$ws.Range("J58").Value = 1200
$ws.Range("L58").Value = 3600
$ws.Range("N58").Value = -3900

# Row 98
$ws.Range("H98").Value = 1216.909
$ws.Range("I98").Value = 481
$ws.Range("K98").Value = 481
$ws.Range("M98").Value = 1017

# Row 122
$ws.Range("H122").Value = 1216.909
$ws.Range("I122").Value = 481
$ws.Range("K122").Value = 1443
$ws.Range("M122").Value = 1007

# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()  # was -60600

# Row 138
$ws.Range("H138").Value = 2393.7666
$ws.Range("I138").Value = 2922.7693
$ws.Range("J138").Value = 2304.4546
$ws.Range("K138").Value = 8768.3079
$ws.Range("L138").Value = 6913.3638
$ws.Range("M138").Value = -3628.3079
$ws.Range("N138").Value = -17193.3638

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 419655.6
$ws.Range("I32").Value = 492835.25
$ws.Range("J32").Value = 17167.5
$ws.Range("K32").Value = 492835.25
$ws.Range("L32").Value = 17167.5
$ws.Range("M32").Value = -492548.25
$ws.Range("N32").Value = -17741.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 97
$ws.Range("H97").Value = 32000
$ws.Range("I97").Value = 20000
$ws.Range("J97").Value = 38000
$ws.Range("K97").Value = 20000
$ws.Range("L97").Value = 38000
$ws.Range("M97").Value = -19009
$ws.Range("N97").Value = -39982

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 270
$ws.Range("I22").Value = 118
$ws.Range("J22").Value = 650
$ws.Range("K22").Value = 118
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = 232
$ws.Range("N22").Value = -1350

# Row 31
$ws.Range("H31").Value = 5640.592
$ws.Range("I31").Value = 1905.625
$ws.Range("J31").Value = 9226.16
$ws.Range("K31").Value = 1905.625
$ws.Range("L31").Value = 9226.16
$ws.Range("M31").Value = -1610.625
$ws.Range("N31").Value = -9816.16

# Row 34
$ws.Range("H34").Value = 5640.592
$ws.Range("I34").Value = 1905.625
$ws.Range("J34").Value = 9226.16
$ws.Range("K34").Value = 1905.625
$ws.Range("L34").Value = 9226.16
$ws.Range("M34").Value = -1703.625
$ws.Range("N34").Value = -9630.16

# Row 99
$ws.Range("H99").Value = 1990.625
$ws.Range("I99").Value = 1850
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1850
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -352
$ws.Range("N99").Value = -4996

# Row 122
$ws.Range("H122").Value = 1839.3043
$ws.Range("I122").Value = 1425
$ws.Range("J122").Value = 1878.762
$ws.Range("K122").Value = 4275
$ws.Range("L122").Value = 5636.286
$ws.Range("M122").Value = -1825
$ws.Range("N122").Value = -10536.286

# Row 126
$ws.Range("H126").Value = 1990.625
$ws.Range("I126").Value = 1850
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5550
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3080
$ws.Range("N126").Value = -10940

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 491.7
$ws.Range("I5").Value = 435.22223
$ws.Range("K5").Value = 1305.66669
$ws.Range("M5").Value = -1193.66669

# Row 68
$ws.Range("H68").Value = 708.5217
$ws.Range("I68").Value = 816
$ws.Range("J68").Value = 625.8461
$ws.Range("K68").Value = 2448
$ws.Range("L68").Value = 1877.5383
$ws.Range("M68").Value = -1637
$ws.Range("N68").Value = -3499.5383

# Row 71
$ws.Range("H71").Value = 708.5217
$ws.Range("I71").Value = 816
$ws.Range("J71").Value = 625.8461
$ws.Range("K71").Value = 7344
$ws.Range("L71").Value = 5632.6149
$ws.Range("M71").Value = -3288
$ws.Range("N71").Value = -13744.6149

# Row 110
$ws.Range("H110").Value = 11446.314
$ws.Range("J110").Value = 12112.3125
$ws.Range("L110").Value = 36336.9375
$ws.Range("N110").Value = -44516.9375

# Row 122
$ws.Range("H122").Value = 8428.691999999999
$ws.Range("I122").Value = 341.42856
$ws.Range("J122").Value = 17863.834
$ws.Range("K122").Value = 3072.85704
$ws.Range("L122").Value = 160774.506
$ws.Range("M122").Value = -622.8570399999999
$ws.Range("N122").Value = -165674.506

# Row 135
$ws.Range("H135").Value = 491.7
$ws.Range("I135").Value = 435.22223
$ws.Range("K135").Value = 3917.00007
$ws.Range("M135").Value = -1382.00007

# Row 139
$ws.Range("H139").Value = 2629.276
$ws.Range("J139").Value = 2856.8572
$ws.Range("L139").Value = 8570.571599999999
$ws.Range("N139").Value = -18850.5716

# Row 141
$ws.Range("H141").Value = 7191.1763
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()  # was 1013.3333

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 4401.1665
$ws.Range("I122").Value = 3007
$ws.Range("J122").Value = 4680
$ws.Range("K122").Value = 9021
$ws.Range("L122").Value = 14040
$ws.Range("M122").Value = -6571
$ws.Range("N122").Value = -18940

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1446.4615
$ws.Range("I7").Value = 1149.9166
$ws.Range("J7").Value = 5005
$ws.Range("K7").Value = 1149.9166
$ws.Range("L7").Value = 5005
$ws.Range("M7").Value = -1037.9166
$ws.Range("N7").Value = -5229

# Row 40
$ws.Range("H40").Value = 502002
$ws.Range("I40").Value = 502002
$ws.Range("K40").Value = 502002
$ws.Range("M40").Value = -501866

# Row 46
$ws.Range("H46").Value = 3871.4285
$ws.Range("I46").Value = 1120
$ws.Range("K46").Value = 1120
$ws.Range("M46").Value = -932

# Row 122
$ws.Range("H122").Value = 3125.4634
$ws.Range("I122").Value = 3035.1428
$ws.Range("J122").Value = 3320
$ws.Range("K122").Value = 9105.428400000001
$ws.Range("L122").Value = 9960
$ws.Range("M122").Value = -6655.428400000001
$ws.Range("N122").Value = -14860

# Row 126
$ws.Range("H126").Value = 1446.4615
$ws.Range("I126").Value = 1149.9166
$ws.Range("J126").Value = 5005
$ws.Range("K126").Value = 3449.7498
$ws.Range("L126").Value = 15015
$ws.Range("M126").Value = -979.7498000000001
$ws.Range("N126").Value = -19955

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4407.769
$ws.Range("I81").Value = 4283.875
$ws.Range("J81").Value = 4606
$ws.Range("K81").Value = 8567.75
$ws.Range("L81").Value = 9212
$ws.Range("M81").Value = -7506.75
$ws.Range("N81").Value = -11334

# Row 84
$ws.Range("H84").Value = 4407.769
$ws.Range("I84").Value = 4283.875
$ws.Range("J84").Value = 4606
$ws.Range("K84").Value = 42838.75
$ws.Range("L84").Value = 46060
$ws.Range("M84").Value = -37534.75
$ws.Range("N84").Value = -56668

# Row 122
$ws.Range("H122").Value = 3840
$ws.Range("I122").Value = 3333.3333
$ws.Range("J122").Value = 4220
$ws.Range("K122").Value = 9999.999899999999
$ws.Range("L122").Value = 12660
$ws.Range("M122").Value = -7549.999899999999
$ws.Range("N122").Value = -17560

# Row 126
$ws.Range("H126").Value = 1994.5555
$ws.Range("I126").Value = 1749.875
$ws.Range("J126").Value = 2190.3
$ws.Range("K126").Value = 5249.625
$ws.Range("L126").Value = 6570.900000000001
$ws.Range("M126").Value = -2779.625
$ws.Range("N126").Value = -11510.9
